# ---------------------------------------------------------------------------
# Helper: replace the contents of a Range with a WordprocessingML body
# fragment via Range.InsertXML (InsertXML replaces the exact range it is
# called on, so the caller must pass the precise Range to overwrite).
# ---------------------------------------------------------------------------
function Set-RangeXml($rng, $innerBodyXml) {
    $full = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
            '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
            '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
            '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
            $innerBodyXml +
            '</w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($full)
}

$RSQUO = [char]0x2019

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Delete the standalone "General Compliance Checklist" paragraph that sits
#    right before the first table.
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("General Compliance Checklist", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para = $r.Paragraphs(1)
$para.Range.Delete()

Write-Host "step1 done"

# ---------------------------------------------------------------------------
# 2. "Proprietary Software" heading: drop its leading lastRenderedPageBreak.
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Proprietary Software", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$inner = '<w:body><w:p><w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr>' +
         '<w:t>Proprietary Software</w:t></w:r></w:p></w:body>'
Set-RangeXml $r $inner

Write-Host "step2 done"

# ---------------------------------------------------------------------------
# 3. "Describe the target's proprietary product or products": add a leading
#    lastRenderedPageBreak.
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Describe the target" + $RSQUO + "s proprietary product or products", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$inner = '<w:body><w:p><w:r><w:rPr><w:bCs/><w:lang w:val="en-US"/></w:rPr>' +
         '<w:lastRenderedPageBreak/><w:t>Describe the target' + $RSQUO + 's proprietary product or products</w:t></w:r></w:p></w:body>'
Set-RangeXml $r $inner

Write-Host "step3 done"
